$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 38 - this pushes the existing rows
# 38..85 down to 39..86 (and the sheet dimension grows to T86).
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with a new market observation,
# following the same constant-column pattern as every other data row.
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = 44483
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100108
$ws.Range("H38").Value = "Tropicales y subtropicales"
$ws.Range("I38").Value = 100108002
$ws.Range("J38").Value = "Mango"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 120
$ws.Range("N38").Value = 7500
$ws.Range("O38").Value = 8000
$ws.Range("P38").Value = 7750
$ws.Range("Q38").Value = "$/bandeja 4 kilos"
$ws.Range("R38").Value = "Perú"
$ws.Range("S38").Value = 1938
$ws.Range("T38").Value = 4
